# Doing Updates for Financials
# Updates yearly financial figures (Income Statement, Balance Sheet, Cash Flow Statement)
# for NAKD on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 89000
$ws.Range("E8").Value = 65200
$ws.Range("F8").Value = 102300
$ws.Range("G8").Value = 94100
$ws.Range("D9").Value = 59300
$ws.Range("E9").Value = 38700
$ws.Range("F9").Value = 56600
$ws.Range("G9").Value = 53500
$ws.Range("D10").Value = 29800
$ws.Range("E10").Value = 26500
$ws.Range("F10").Value = 45700
$ws.Range("G10").Value = 40500
$ws.Range("F14").Value = 3000
$ws.Range("G14").Value = 8300
$ws.Range("D17").Value = 114400
$ws.Range("E17").Value = 75500
$ws.Range("F17").Value = 112600
$ws.Range("G17").Value = 103800
$ws.Range("D18").Value = -25400
$ws.Range("E18").Value = -10200
$ws.Range("F18").Value = -10300
$ws.Range("G18").Value = -9700
$ws.Range("D21").Value = -23400
$ws.Range("E21").Value = -9000
$ws.Range("F21").Value = -7900
$ws.Range("G21").Value = -5900
$ws.Range("D23").Value = -25400
$ws.Range("E23").Value = -10200
$ws.Range("F23").Value = -10300
$ws.Range("G23").Value = -9700
$ws.Range("F24").Value = 3800
$ws.Range("G24").Value = -900
$ws.Range("D26").Value = -25500
$ws.Range("E26").Value = -10800
$ws.Range("F26").Value = -14100
$ws.Range("G26").Value = -8900
$ws.Range("D27").Value = -25500
$ws.Range("E27").Value = -10800
$ws.Range("F27").Value = -14100
$ws.Range("G27").Value = -8900
$ws.Range("D33").Value = -25500
$ws.Range("E33").Value = -10800
$ws.Range("F33").Value = -14100
$ws.Range("G33").Value = -8900
$ws.Range("D35").Value = -25500
$ws.Range("E35").Value = -10800
$ws.Range("F35").Value = -14100
$ws.Range("G35").Value = -8900
$ws.Range("D41").Value = 7300
$ws.Range("D43").Value = 19100
$ws.Range("E43").Value = 27800
$ws.Range("F43").Value = 22100
$ws.Range("G43").Value = 13600
$ws.Range("D44").Value = 21300
$ws.Range("E44").Value = 25700
$ws.Range("F44").Value = 25700
$ws.Range("G44").Value = 31500
$ws.Range("G45").Value = 1600
$ws.Range("D46").Value = 47700
$ws.Range("E46").Value = 55300
$ws.Range("F46").Value = 50700
$ws.Range("G46").Value = 47400
$ws.Range("E48").Value = 3400
$ws.Range("F48").Value = 4200
$ws.Range("D49").Value = 8800
$ws.Range("E49").Value = 9900
$ws.Range("F49").Value = 9900
$ws.Range("G49").Value = 11800
$ws.Range("G52").Value = 3800
$ws.Range("D54").Value = 59700
$ws.Range("E54").Value = 68600
$ws.Range("F54").Value = 64800
$ws.Range("G54").Value = 67600
$ws.Range("D57").Value = 15500
$ws.Range("E57").Value = 14300
$ws.Range("F57").Value = 13500
$ws.Range("G57").Value = 18600
$ws.Range("D58").Value = 35300
$ws.Range("E58").Value = 46700
$ws.Range("F58").Value = 41700
$ws.Range("G58").Value = 38100
$ws.Range("D59").Value = 10900
$ws.Range("E59").Value = 12200
$ws.Range("F59").Value = 9000
$ws.Range("G59").Value = 7000
$ws.Range("D60").Value = 61700
$ws.Range("E60").Value = 73200
$ws.Range("F60").Value = 64200
$ws.Range("G60").Value = 63700
$ws.Range("F61").Value = 10800
$ws.Range("G62").Value = 2000
$ws.Range("D66").Value = 63600
$ws.Range("E66").Value = 74700
$ws.Range("F66").Value = 76900
$ws.Range("G66").Value = 65700
$ws.Range("D72").Value = -50400
$ws.Range("E72").Value = -25100
$ws.Range("F72").Value = -14200
$ws.Range("D76").Value = -3900
$ws.Range("E76").Value = -6100
$ws.Range("F76").Value = -12100
$ws.Range("D81").Value = -25500
$ws.Range("E81").Value = -10800
$ws.Range("F81").Value = -14100
$ws.Range("G81").Value = -8900
$ws.Range("D83").Value = 2100
$ws.Range("F83").Value = 2400
$ws.Range("G83").Value = 3900
$ws.Range("D89").Value = -2800
$ws.Range("E89").Value = -9200
$ws.Range("G89").Value = -11700
$ws.Range("E91").Value = -500
$ws.Range("F91").Value = -1800
$ws.Range("G91").Value = -3200
$ws.Range("D94").Value = -1600
$ws.Range("F94").Value = -2200
$ws.Range("D100").Value = 9800
$ws.Range("E100").Value = 8900
$ws.Range("F100").Value = 7600
$ws.Range("G100").Value = 13900
$ws.Range("D102").Value = 5500
